$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 458.07144
$ws.Cells.Item(17, 10).Value = 458.07144
$ws.Cells.Item(17, 12).Value = 1374.21432
$ws.Cells.Item(17, 14).Value = -1710.21432

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 50906.5
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 13).ClearContents()
$ws.Cells.Item(32, 8).Value = 3553.8713
$ws.Cells.Item(32, 9).Value = 3278.4915
$ws.Cells.Item(32, 10).Value = 5030.909
$ws.Cells.Item(32, 11).Value = 3278.4915
$ws.Cells.Item(32, 12).Value = 5030.909
$ws.Cells.Item(32, 13).Value = -2991.4915
$ws.Cells.Item(32, 14).Value = -5604.909
$ws.Cells.Item(61, 8).Value = 1089.841
$ws.Cells.Item(61, 9).Value = 950.8570999999999
$ws.Cells.Item(61, 11).Value = 950.8570999999999
$ws.Cells.Item(61, 13).Value = -738.8570999999999
$ws.Cells.Item(74, 8).Value = 1663.1333
$ws.Cells.Item(74, 9).Value = 991.9
$ws.Cells.Item(74, 11).Value = 991.9
$ws.Cells.Item(74, 13).Value = -117.9
$ws.Cells.Item(77, 8).Value = 1663.1333
$ws.Cells.Item(77, 9).Value = 991.9
$ws.Cells.Item(77, 11).Value = 4959.5
$ws.Cells.Item(77, 13).Value = -591.5
$ws.Cells.Item(110, 8).Value = 2325.111
$ws.Cells.Item(110, 9).Value = 1633.3334
$ws.Cells.Item(110, 10).Value = 2671
$ws.Cells.Item(110, 11).Value = 1633.3334
$ws.Cells.Item(110, 12).Value = 2671
$ws.Cells.Item(110, 13).Value = 411.6666
$ws.Cells.Item(110, 14).Value = -6761
$ws.Cells.Item(116, 8).Value = 50906.5
$ws.Cells.Item(116, 9).Value = 0
$ws.Cells.Item(116, 11).Value = 0
$ws.Cells.Item(116, 13).ClearContents()
$ws.Cells.Item(132, 8).Value = 1730.1923
$ws.Cells.Item(132, 9).Value = 1408.5
$ws.Cells.Item(132, 11).Value = 4225.5
$ws.Cells.Item(132, 13).Value = -1695.5
$ws.Cells.Item(136, 8).Value = 1089.841
$ws.Cells.Item(136, 9).Value = 950.8570999999999
$ws.Cells.Item(136, 11).Value = 2852.5713
$ws.Cells.Item(136, 13).Value = -302.5712999999996
$ws.Cells.Item(139, 8).Value = 39175
$ws.Cells.Item(139, 10).Value = 39175
$ws.Cells.Item(139, 12).Value = 39175
$ws.Cells.Item(139, 14).Value = -49455

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 50906.5
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 13).ClearContents()
$ws.Cells.Item(86, 8).Value = 4515.76
$ws.Cells.Item(86, 9).Value = 4540.6816
$ws.Cells.Item(86, 11).Value = 4540.6816
$ws.Cells.Item(86, 13).Value = -3417.6816
$ws.Cells.Item(89, 8).Value = 4515.76
$ws.Cells.Item(89, 9).Value = 4540.6816
$ws.Cells.Item(89, 11).Value = 22703.408
$ws.Cells.Item(89, 13).Value = -17087.408
$ws.Cells.Item(105, 8).Value = 55557956
$ws.Cells.Item(105, 9).Value = 58825984
$ws.Cells.Item(105, 10).Value = 1500
$ws.Cells.Item(105, 11).Value = 58825984
$ws.Cells.Item(105, 12).Value = 1500
$ws.Cells.Item(105, 13).Value = -58824237
$ws.Cells.Item(105, 14).Value = -4994
$ws.Cells.Item(134, 8).Value = 4642.4414
$ws.Cells.Item(134, 9).Value = 1088
$ws.Cells.Item(134, 10).Value = 14515.889
$ws.Cells.Item(134, 11).Value = 3264
$ws.Cells.Item(134, 12).Value = 43547.667
$ws.Cells.Item(134, 13).Value = -729
$ws.Cells.Item(134, 14).Value = -48617.667
$ws.Cells.Item(138, 8).Value = 41256.668
$ws.Cells.Item(138, 10).Value = 41256.668
$ws.Cells.Item(138, 12).Value = 41256.668
$ws.Cells.Item(138, 14).Value = -51536.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1087.8182
$ws.Cells.Item(58, 9).Value = 858.25
$ws.Cells.Item(58, 10).Value = 1700
$ws.Cells.Item(58, 11).Value = 858.25
$ws.Cells.Item(58, 12).Value = 1700
$ws.Cells.Item(58, 13).Value = -655.25
$ws.Cells.Item(58, 14).Value = -2106
$ws.Cells.Item(99, 8).Value = 2393885.2
$ws.Cells.Item(99, 9).Value = 2925415.5
$ws.Cells.Item(99, 10).Value = 1999.5
$ws.Cells.Item(99, 11).Value = 2925415.5
$ws.Cells.Item(99, 12).Value = 1999.5
$ws.Cells.Item(99, 13).Value = -2923917.5
$ws.Cells.Item(99, 14).Value = -4995.5
$ws.Cells.Item(126, 8).Value = 2393885.2
$ws.Cells.Item(126, 9).Value = 2925415.5
$ws.Cells.Item(126, 10).Value = 1999.5
$ws.Cells.Item(126, 11).Value = 8776246.5
$ws.Cells.Item(126, 12).Value = 5998.5
$ws.Cells.Item(126, 13).Value = -8773776.5
$ws.Cells.Item(126, 14).Value = -10938.5
$ws.Cells.Item(132, 8).Value = 2646.611
$ws.Cells.Item(132, 9).Value = 2049.3076
$ws.Cells.Item(132, 11).Value = 6147.9228
$ws.Cells.Item(132, 13).Value = -3617.9228
$ws.Cells.Item(136, 8).Value = 1087.8182
$ws.Cells.Item(136, 9).Value = 858.25
$ws.Cells.Item(136, 10).Value = 1700
$ws.Cells.Item(136, 11).Value = 2574.75
$ws.Cells.Item(136, 12).Value = 5100
$ws.Cells.Item(136, 13).Value = -24.75
$ws.Cells.Item(136, 14).Value = -10200

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(98, 8).Value = 466.33334
$ws.Cells.Item(98, 10).Value = 466.33334
$ws.Cells.Item(98, 12).Value = 1399.00002
$ws.Cells.Item(98, 14).Value = -4395.000019999999
$ws.Cells.Item(103, 8).Value = 1506
$ws.Cells.Item(103, 9).Value = 181.66667
$ws.Cells.Item(103, 10).Value = 1947.4445
$ws.Cells.Item(103, 11).Value = 545.00001
$ws.Cells.Item(103, 12).Value = 5842.333500000001
$ws.Cells.Item(103, 13).Value = 333.99999
$ws.Cells.Item(103, 14).Value = -7600.333500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1746.75
$ws.Cells.Item(102, 9).Value = 1746.75
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 1746.75
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 14).ClearContents()
$ws.Cells.Item(102, 13).Value = -124.75
$ws.Cells.Item(107, 8).Value = 687442.8
$ws.Cells.Item(107, 9).Value = 1069050
$ws.Cells.Item(107, 11).Value = 1069050
$ws.Cells.Item(107, 13).Value = -1067130
$ws.Cells.Item(113, 8).Value = 1608
$ws.Cells.Item(113, 9).Value = 1600
$ws.Cells.Item(113, 10).Value = 1614.6666
$ws.Cells.Item(113, 11).Value = 1600
$ws.Cells.Item(113, 12).Value = 1614.6666
$ws.Cells.Item(113, 13).Value = 570
$ws.Cells.Item(113, 14).Value = -5954.6666
$ws.Cells.Item(132, 8).Value = 2348.0334
$ws.Cells.Item(132, 9).Value = 1739.4615
$ws.Cells.Item(132, 11).Value = 5218.3845
$ws.Cells.Item(132, 13).Value = -2688.3845

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2077.7778
$ws.Cells.Item(7, 9).Value = 1840
$ws.Cells.Item(7, 10).Value = 2375
$ws.Cells.Item(7, 11).Value = 1840
$ws.Cells.Item(7, 12).Value = 2375
$ws.Cells.Item(7, 13).Value = -1728
$ws.Cells.Item(7, 14).Value = -2599
$ws.Cells.Item(16, 8).Value = 630.3333
$ws.Cells.Item(16, 9).Value = 630.3333
$ws.Cells.Item(16, 11).Value = 630.3333
$ws.Cells.Item(16, 13).Value = -460.3333
$ws.Cells.Item(22, 8).Value = 1948.7142
$ws.Cells.Item(22, 10).Value = 1948.7142
$ws.Cells.Item(22, 12).Value = 1948.7142
$ws.Cells.Item(22, 14).Value = -2538.7142
$ws.Cells.Item(27, 8).Value = 1948.7142
$ws.Cells.Item(27, 10).Value = 1948.7142
$ws.Cells.Item(27, 12).Value = 1948.7142
$ws.Cells.Item(27, 14).Value = -2162.7142
$ws.Cells.Item(46, 8).Value = 2599.8
$ws.Cells.Item(46, 9).Value = 1000
$ws.Cells.Item(46, 10).Value = 3666.3333
$ws.Cells.Item(46, 11).Value = 1000
$ws.Cells.Item(46, 12).Value = 3666.3333
$ws.Cells.Item(46, 13).Value = -812
$ws.Cells.Item(46, 14).Value = -4042.3333
$ws.Cells.Item(61, 8).Value = 1232.4736
$ws.Cells.Item(61, 9).Value = 945.1818
$ws.Cells.Item(61, 10).Value = 1627.5
$ws.Cells.Item(61, 11).Value = 945.1818
$ws.Cells.Item(61, 12).Value = 1627.5
$ws.Cells.Item(61, 13).Value = -743.1818
$ws.Cells.Item(61, 14).Value = -2031.5
$ws.Cells.Item(113, 8).Value = 1232.4736
$ws.Cells.Item(113, 9).Value = 945.1818
$ws.Cells.Item(113, 10).Value = 1627.5
$ws.Cells.Item(113, 11).Value = 945.1818
$ws.Cells.Item(113, 12).Value = 1627.5
$ws.Cells.Item(113, 13).Value = 1224.8182
$ws.Cells.Item(113, 14).Value = -5967.5
$ws.Cells.Item(122, 8).Value = 11807931
$ws.Cells.Item(122, 9).Value = 25759720
$ws.Cells.Item(122, 10).Value = 2570.923
$ws.Cells.Item(122, 11).Value = 77279160
$ws.Cells.Item(122, 12).Value = 7712.768999999999
$ws.Cells.Item(122, 13).Value = -77276710
$ws.Cells.Item(122, 14).Value = -12612.769
$ws.Cells.Item(126, 8).Value = 2077.7778
$ws.Cells.Item(126, 9).Value = 1840
$ws.Cells.Item(126, 10).Value = 2375
$ws.Cells.Item(126, 11).Value = 5520
$ws.Cells.Item(126, 12).Value = 7125
$ws.Cells.Item(126, 13).Value = -3050
$ws.Cells.Item(126, 14).Value = -12065
$ws.Cells.Item(132, 8).Value = 23277.936
$ws.Cells.Item(132, 9).Value = 1016.62067
$ws.Cells.Item(132, 11).Value = 3049.86201
$ws.Cells.Item(132, 13).Value = -519.8620099999998
$ws.Cells.Item(136, 8).Value = 1381.0588
$ws.Cells.Item(136, 9).Value = 1385.5555
$ws.Cells.Item(136, 10).Value = 1376
$ws.Cells.Item(136, 11).Value = 4156.666499999999
$ws.Cells.Item(136, 12).Value = 4128
$ws.Cells.Item(136, 13).Value = -1606.666499999999
$ws.Cells.Item(136, 14).Value = -9228

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 0
$ws.Cells.Item(46, 10).Value = 0
$ws.Cells.Item(46, 14).ClearContents()
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(69, 8).Value = 0
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 14).ClearContents()
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(72, 8).Value = 0
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 14).ClearContents()
$ws.Cells.Item(72, 12).Value = 0
$ws.Cells.Item(100, 8).Value = 534.6667
$ws.Cells.Item(100, 9).Value = 534.6667
$ws.Cells.Item(100, 11).Value = 1069.3334
$ws.Cells.Item(100, 13).Value = -528.3334
$ws.Cells.Item(107, 8).Value = 426.2
$ws.Cells.Item(107, 9).Value = 430.25
$ws.Cells.Item(107, 10).Value = 410
$ws.Cells.Item(107, 11).Value = 1290.75
$ws.Cells.Item(107, 12).Value = 1230
$ws.Cells.Item(107, 13).Value = 629.25
$ws.Cells.Item(107, 14).Value = -5070
$ws.Cells.Item(126, 8).Value = 71430050
$ws.Cells.Item(126, 9).Value = 90910540
$ws.Cells.Item(126, 10).Value = 1566.6666
$ws.Cells.Item(126, 11).Value = 272731620
$ws.Cells.Item(126, 12).Value = 4699.9998
$ws.Cells.Item(126, 13).Value = -272729150
$ws.Cells.Item(126, 14).Value = -9639.9998
$ws.Cells.Item(132, 8).Value = 1849.1666
$ws.Cells.Item(132, 9).Value = 1051.1666
$ws.Cells.Item(132, 11).Value = 3153.4998
$ws.Cells.Item(132, 13).Value = -623.4998000000001
$ws.Cells.Item(134, 8).Value = 0
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 14).ClearContents()
$ws.Cells.Item(134, 12).Value = 0
